# Apply Lrpap1-Vldlr LR-pair data update (Natmi following Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @("ECs","Lrpap1","Vldlr","ECs",3,1,6.660188000000001,19.980564,0.1500148400131262,0.1500148400131261,1,0.3333333333333333,0.1112246666666667,0.333674,0.008007360146089436,0.008007360146089434,0.7407771902373335,6.666994712136001,0.001201222851243089,0.001201222851243089),
    @("ECs","Lrpap1","Vldlr","FAPs",3,1,6.660188000000001,19.980564,0.1500148400131262,0.1500148400131261,3,1,11.825228,35.475684,0.851329675722,0.8513296757219999,78.75824162286402,708.8241746057761,0.1277120851018624,0.1277120851018624),
    @("ECs","Lrpap1","Vldlr","M2",3,1,6.660188000000001,19.980564,0.1500148400131262,0.1500148400131261,1,0.3333333333333333,0.03793733333333333,0.113812,0.002731209722503793,0.002731209722503793,0.2526697722186667,2.274027949968,0.0004097219895637013,0.0004097219895637012),
    @("ECs","Lrpap1","Vldlr","sCs",3,1,6.660188000000001,19.980564,0.1500148400131262,0.1500148400131261,3,1,1.915914,5.747742000000001,0.1379317544094067,0.1379317544094067,12.760347431832,114.843126886488,0.02069181007045696,0.02069181007045696),
    @("FAPs","Lrpap1","Vldlr","ECs",3,1,16.410331,49.230993,0.3696281815959916,0.3696281815959916,1,0.3333333333333333,0.1112246666666667,0.333674,0.008007360146089436,0.008007360146089434,1.825233595364667,16.427102358282,0.002959745970183252,0.002959745970183252),
    @("FAPs","Lrpap1","Vldlr","FAPs",3,1,16.410331,49.230993,0.3696281815959916,0.3696281815959916,3,1,11.825228,35.475684,0.851329675722,0.8513296757219999,194.055905630468,1746.503150674212,0.3146754399758281,0.3146754399758281),
    @("FAPs","Lrpap1","Vldlr","M2",3,1,16.410331,49.230993,0.3696281815959916,0.3696281815959916,1,0.3333333333333333,0.03793733333333333,0.113812,0.002731209722503793,0.002731209722503793,0.6225641972573333,5.603077775316,0.00100953208328637,0.00100953208328637),
    @("FAPs","Lrpap1","Vldlr","sCs",3,1,16.410331,49.230993,0.3696281815959916,0.3696281815959916,3,1,1.915914,5.747742000000001,0.1379317544094067,0.1379317544094067,31.440782907534,282.967046167806,0.05098346356669391,0.05098346356669391),
    @("M2","Lrpap1","Vldlr","ECs",3,1,14.478895,43.436685,0.3261242951387937,0.3261242951387937,1,0.3333333333333333,0.1112246666666667,0.333674,0.008007360146089436,0.008007360146089434,1.610410270076667,14.49369243069,0.002611394683565885,0.002611394683565884),
    @("M2","Lrpap1","Vldlr","FAPs",3,1,14.478895,43.436685,0.3261242951387937,0.3261242951387937,3,1,11.825228,35.475684,0.851329675722,0.8513296757219999,171.21623456306,1540.94611106754,0.2776392904255751,0.277639290425575),
    @("M2","Lrpap1","Vldlr","M2",3,1,14.478895,43.436685,0.3261242951387937,0.3261242951387937,1,0.3333333333333333,0.03793733333333333,0.113812,0.002731209722503793,0.002731209722503793,0.5492906659133333,4.94361599322,0.0008907138456277699,0.0008907138456277696),
    @("M2","Lrpap1","Vldlr","sCs",3,1,14.478895,43.436685,0.3261242951387937,0.3261242951387937,3,1,1.915914,5.747742000000001,0.1379317544094067,0.1379317544094067,27.74031763503,249.66285871527,0.04498289618402497,0.04498289618402496),
    @("sCs","Lrpap1","Vldlr","ECs",3,1,6.847447,20.542341,0.1542326832520885,0.1542326832520885,1,0.3333333333333333,0.1112246666666667,0.333674,0.008007360146089436,0.008007360146089434,0.7616050100926668,6.854445090834001,0.001234996641097209,0.001234996641097208),
    @("sCs","Lrpap1","Vldlr","FAPs",3,1,6.847447,20.542341,0.1542326832520885,0.1542326832520885,3,1,11.825228,35.475684,0.851329675722,0.8513296757219999,80.972621992916,728.7535979362441,0.1313028602187344,0.1313028602187344),
    @("sCs","Lrpap1","Vldlr","M2",3,1,6.847447,20.542341,0.1542326832520885,0.1542326832520885,1,0.3333333333333333,0.03793733333333333,0.113812,0.002731209722503793,0.002731209722503793,0.2597738793213333,2.337964913892,0.0004212418040259519,0.0004212418040259519),
    @("sCs","Lrpap1","Vldlr","sCs",3,1,6.847447,20.542341,0.1542326832520885,0.1542326832520885,3,1,1.915914,5.747742000000001,0.1379317544094067,0.1379317544094067,13.119119571558,118.072076144022,0.02127358458823088,0.02127358458823088)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $rowVals = $rowsData[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $rowVals[$c - 1]
    }
}
